# Update "想去人数" (want-to-go count) figures on the "展览" and "全部类型" sheets
# to match the newly scraped totals.

$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibitions)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value = 2050   # was 2047 - 南宁·草莓动漫节
$ws1.Range("F4").Value = 853    # was 849  - 南宁·第一届ANE·DACG动漫嘉年华
$ws1.Range("F5").Value = 1145   # was 1133 - 南宁·2024三月三国潮动漫节（良牙春典）

# Sheet "全部类型" (All types)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F3").Value = 2050   # was 2047 - 南宁·草莓动漫节
$ws4.Range("F6").Value = 853    # was 849  - 南宁·第一届ANE·DACG动漫嘉年华
$ws4.Range("F7").Value = 1145   # was 1133 - 南宁·2024三月三国潮动漫节（良牙春典）
